$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.29025
$ws.Range("H2").Value = 6.870749999999999
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 5.855348
$ws.Range("N2").Value = 17.566044
$ws.Range("O2").Value = 0.05092948808292105
$ws.Range("P2").Value = 0.05092948808292105
$ws.Range("Q2").Value = 13.410210757
$ws.Range("R2").Value = 120.691896813
$ws.Range("S2").Value = 0.05092948808292105
$ws.Range("T2").Value = 0.05092948808292105

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.29025
$ws.Range("H3").Value = 6.870749999999999
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 62.99699166666667
$ws.Range("N3").Value = 188.990975
$ws.Range("O3").Value = 0.547944295769846
$ws.Range("P3").Value = 0.547944295769846
$ws.Range("Q3").Value = 144.2788601645833
$ws.Range("R3").Value = 1298.50974148125
$ws.Range("S3").Value = 0.547944295769846
$ws.Range("T3").Value = 0.547944295769846

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.29025
$ws.Range("H4").Value = 6.870749999999999
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.2640463333333333
$ws.Range("N4").Value = 0.792139
$ws.Range("O4").Value = 0.002296660179179615
$ws.Range("P4").Value = 0.002296660179179615
$ws.Range("Q4").Value = 0.6047321149166667
$ws.Range("R4").Value = 5.44258903425
$ws.Range("S4").Value = 0.002296660179179615
$ws.Range("T4").Value = 0.002296660179179615

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.29025
$ws.Range("H5").Value = 6.870749999999999
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.226170666666667
$ws.Range("N5").Value = 3.678512
$ws.Range("O5").Value = 0.01066516360011862
$ws.Range("P5").Value = 0.01066516360011862
$ws.Range("Q5").Value = 2.808237369333333
$ws.Range("R5").Value = 25.274136324
$ws.Range("S5").Value = 0.01066516360011862
$ws.Range("T5").Value = 0.01066516360011862

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.29025
$ws.Range("H6").Value = 6.870749999999999
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.763133
$ws.Range("N6").Value = 2.289399
$ws.Range("O6").Value = 0.006637687978440185
$ws.Range("P6").Value = 0.006637687978440185
$ws.Range("Q6").Value = 1.74776535325
$ws.Range("R6").Value = 15.72988817925
$ws.Range("S6").Value = 0.006637687978440185
$ws.Range("T6").Value = 0.006637687978440185

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.29025
$ws.Range("H7").Value = 6.870749999999999
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 43.864011
$ws.Range("N7").Value = 131.592033
$ws.Range("O7").Value = 0.3815267043894945
$ws.Range("P7").Value = 0.3815267043894945
$ws.Range("Q7").Value = 100.45955119275
$ws.Range("R7").Value = 904.13596073475
$ws.Range("S7").Value = 0.3815267043894945
$ws.Range("T7").Value = 0.3815267043894945

